# Task Assignments.xlsx - "Add files via upload" re-upload of an earlier
# save of the workbook: Build phase is finished (column D -> 100%), the
# Testing phase has now also been run and completed (columns F, G, H),
# the rows that previously had no Build dates got them, and the "Names"
# column (E) plus the new Testing "Names" column (I) are filled in with
# the same assignee, "Vĩnh Kha", instead of the stray "Thế Anh"/blank
# entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Fill in the previously-blank Build start/finish dates (rows 7-11).
# ---------------------------------------------------------------------
$ws.Range("B7").Value = 44835
$ws.Range("C7").Value = 44835

$ws.Range("B8").Value = 44843
$ws.Range("C8").Value = 44844

$ws.Range("B9").Value = 44845
$ws.Range("C9").Value = 44846

$ws.Range("B10").Value = 44847
$ws.Range("C10").Value = 44847

$ws.Range("B11").Value = 44848
$ws.Range("C11").Value = 44848

# Rows 2-6 already carried real Build dates, but some of those cells
# (7, 9 & 10) were stamped with the stray "empty-date" style (no fill) -
# bring every Start/Finish Build cell onto the same look as B2 (the
# plain dated style) so the freshly-populated cells match their
# neighbours.
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("C7").PasteSpecial(-4122)
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Build is done for every task -> % Build = 100%.
# ---------------------------------------------------------------------
$ws.Range("D2:D11").Value = 1

# Give every % Build cell the same "complete" look (green fill) as D2.
$ws.Range("D2").Copy()
$ws.Range("D5:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 3. Testing has now been carried out for every task: Start/Finish
#    Testing = 14 Oct 2022 (44848), % Testing = 100%.
# ---------------------------------------------------------------------
$ws.Range("F2:F11").Value = 44848
$ws.Range("G2:G11").Value = 44848
$ws.Range("H2:H11").Value = 1

# Start/Finish Testing should look like the Build date columns.
$ws.Range("B2").Copy()
$ws.Range("F2:F11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# % Testing should look like the now-green % Build column.
$ws.Range("D2").Copy()
$ws.Range("H2:H11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. Names: both the Build "Names" (E) and Testing "Names" (I) columns
#    are assigned to "Vĩnh Kha" for every task.
# ---------------------------------------------------------------------
$ws.Range("E2:E11").Value = "Vĩnh Kha"
$ws.Range("I2:I11").Value = "Vĩnh Kha"

# ---------------------------------------------------------------------
# 5. Cosmetic touch-ups that came along with the re-save.
# ---------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 24.75
$ws.Range("A2:J11").RowHeight = 15.75

$ws.Columns.Item(1).ColumnWidth = 13.43
$ws.Columns.Item(2).ColumnWidth = 10.75

$ws.Range("I2:I11").Select()
